# I0 and IF added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style of the existing
# header cells (e.g. H1) so the look/feel is consistent.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-33: column I is always 1, column J mirrors column H (IP).
for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $ws.Cells.Item($row, 8).Value2
}
